# Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study
#
# The CasesTab query (cell B2 on the "startup" sheet) incorrectly joined to
# an optional :cohort node and returned a `Cohort` column that isn't part of
# this test's expected output. Remove the optional cohort match/column so the
# query again lines up with the other fixed test queries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.sex IN ['Female']
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $casesQuery

# The workbook was left with the selection sitting on the (now shorter)
# CasesTab query cell instead of scrolled down to the FilesTab row.
$ws.Range("B2").Select()
